$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the install type for RHEL (row 3) from t2.nano to t2.micro
$ws.Range("C3").Value = "t2.micro"

# Move the active cell selection to C4
$ws.Activate()
$ws.Range("C4").Select()
